# FINFLUX-2815 Stabilaizing automation script
# Adds a new "Modify Transaction1" sheet (loan-navigation helper used by the
# automation suite), tweaks the Summary sheet's last "Over Due" figure back
# to 0, and nudges a couple of UI-only bits (selection / first visible tab)
# left over from the authoring session.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------
# 1. Update the Summary sheet: F5 (Over Due total) goes back to 0, and
#    remember the selection the author left the cursor on.
# ---------------------------------------------------------------------
$summary = $sheets.Item("Summary")
$summary.Range("F5").Value = 0
$summary.Range("C9").Select()

# ---------------------------------------------------------------------
# 2. Append a brand-new worksheet at the very end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $sheets.Item($sheets.Count)
$newSheet = $sheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Modify Transaction1"

$newSheet.Range("A1").Value = "OverDueTillDate"
$newSheet.Range("B1").Value = 42064
$newSheet.Range("B1").NumberFormat = "d-mmm-yy"

$newSheet.Range("A2").Value = "clickonsubmit"
$newSheet.Range("B2").Value = "Submit"

$newSheet.Range("A3").Value = "NavigateToLoan"
$newSheet.Range("B3").Value = "navigate"

$newSheet.Range("A1:A3").Font.Bold = $false
$newSheet.Range("A1:A3").Interior.Color = 15921906

$newSheet.Columns("A").ColumnWidth = 15.85546875

$newSheet.Range("D9").Select()

# ---------------------------------------------------------------------
# 3. Re-activate the Summary tab (adding a sheet switches focus to it)
#    and scroll the tab strip so the third tab is first visible, matching
#    the state the workbook was saved in.
# ---------------------------------------------------------------------
$summary.Activate()
$wb.Windows.Item(1).ScrollWorkbookTabs(1, 3)
